$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.004.92"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "2.419.63"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.57%  "
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D16").Value = "61.899.61"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "2.411.29"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "323.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "554.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.83%  "
$ws.Range("D27").Value = "2.536.29"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "0.0" + [string]([char]0x2083) + "0936"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("E31").Value = "  -3.75%  "
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "153.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.992"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "147.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.94%  "
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.592"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("E51").Value = "  +0.57%  "
